# ajuste na figura echo
# Adds a new "LV_mass_z score" column (M) with values for every data row,
# mirroring the formatting already used by the neighbouring H/L columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131
$xlCenter = -4108
$xlPasteValues = -4163

# ---- Header (M1) ----------------------------------------------------
$ws.Range("M1").Value = "LV_mass_z score"
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").Font.Name = "Arial"
$ws.Range("M1").HorizontalAlignment = $xlLeft

# ---- Rows whose text is entered with a leading apostrophe -----------
# (these keep the "quote prefixed" look already used in columns H/L)
$quotedRows = @{
    2 = "-0.91"
    3 = "-1.05"
    4 = "-1.36"
}
foreach ($r in $quotedRows.Keys) {
    $cell = $ws.Range("M" + $r)
    $cell.NumberFormat = "0"
    $cell.HorizontalAlignment = $xlCenter
    $cell.Value = "'" + $quotedRows[$r]
}

# ---- Remaining rows: plain text values, centred, no quote prefix ----
$plainRows = @{
    5  = "0.51"
    6  = "0.49"
    7  = "2.62"
    8  = "1.48"
    9  = "0.33"
    10 = "0.87"
    11 = "1.53"
    12 = "1.76"
    13 = "0.68"
}

$scratch = $ws.Range("P20")
foreach ($r in $plainRows.Keys) {
    $target = $ws.Range("M" + $r)
    $target.NumberFormat = "0"
    $target.HorizontalAlignment = $xlCenter

    $scratch.Value = "'" + $plainRows[$r]
    $scratch.Copy()
    $target.PasteSpecial($xlPasteValues)
}
$excel.CutCopyMode = $false
$scratch.Clear()

# ---- Selection / view state ------------------------------------------
$ws.Range("M2:M13").Select()
